$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 5100
$ws.Range("J48").Value = 5100
$ws.Range("L48").Value = 15300
$ws.Range("N48").Value = -15884
$ws.Range("H56").Value = 5100
$ws.Range("J56").Value = 5100
$ws.Range("L56").Value = 15300
$ws.Range("N56").Value = -16368
$ws.Range("H113").Value = 68637
$ws.Range("I113").Value = 88183.164
$ws.Range("K113").Value = 88183.164
$ws.Range("M113").Value = -84929.164
$ws.Range("H125").Value = 3578.1538
$ws.Range("I125").Value = 3768.75
$ws.Range("J125").Value = 3273.2
$ws.Range("K125").Value = 33918.75
$ws.Range("L125").Value = 29458.8
$ws.Range("M125").Value = -31458.75
$ws.Range("N125").Value = -34378.8
$ws.Range("H135").Value = 1410.8
$ws.Range("I135").Value = 1466.5
$ws.Range("K135").Value = 13198.5
$ws.Range("M135").Value = -10663.5
$ws.Range("H137").Value = 1449.5
$ws.Range("I137").Value = 1449.5
$ws.Range("K137").Value = 4348.5
$ws.Range("M137").Value = -1798.5
$ws.Range("H141").Value = 3642.5217
$ws.Range("I141").Value = 3656.6191
$ws.Range("K141").Value = 10969.8573
$ws.Range("M141").Value = -5789.8573

# --- Worksheet: ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4638.5654
$ws.Range("I32").Value = 4019.862
$ws.Range("J32").Value = 7900.8184
$ws.Range("K32").Value = 4019.862
$ws.Range("L32").Value = 7900.8184
$ws.Range("M32").Value = -3732.862
$ws.Range("N32").Value = -8474.8184
$ws.Range("H47").Value = 45250
$ws.Range("J47").Value = 45250
$ws.Range("L47").Value = 45250
$ws.Range("N47").Value = -46700
$ws.Range("H95").Value = 52333.332
$ws.Range("J95").Value = 52333.332
$ws.Range("L95").Value = 52333.332
$ws.Range("N95").Value = -57825.332
$ws.Range("H97").Value = 2391.6875
$ws.Range("I97").Value = 568.8182
$ws.Range("J97").Value = 6402
$ws.Range("K97").Value = 568.8182
$ws.Range("L97").Value = 6402
$ws.Range("M97").Value = -72.81820000000005
$ws.Range("N97").Value = -7394

# --- Worksheet: BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 62974.25
$ws.Range("J74").Value = 62974.25
$ws.Range("L74").Value = 62974.25
$ws.Range("N74").Value = -64846.25
$ws.Range("H77").Value = 62974.25
$ws.Range("J77").Value = 62974.25
$ws.Range("L77").Value = 188922.75
$ws.Range("N77").Value = -198282.75
$ws.Range("H81").Value = 21107.572
$ws.Range("J81").Value = 21107.572
$ws.Range("L81").Value = 21107.572
$ws.Range("N81").Value = -23229.572
$ws.Range("H84").Value = 21107.572
$ws.Range("J84").Value = 21107.572
$ws.Range("L84").Value = 63322.716
$ws.Range("N84").Value = -73930.716
$ws.Range("H26").Value = 28471
$ws.Range("I26").Value = 28471
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 28471
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -28179
$ws.Range("N26").ClearContents()

# --- Worksheet: CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 115468.5
$ws.Range("J52").Value = 115468.5
$ws.Range("L52").Value = 115468.5
$ws.Range("N52").Value = -116056.5
$ws.Range("H95").Value = 30574
$ws.Range("J95").Value = 30574
$ws.Range("L95").Value = 30574
$ws.Range("N95").Value = -36066
$ws.Range("H96").Value = 1911.5
$ws.Range("J96").Value = 1911.5
$ws.Range("L96").Value = 1911.5
$ws.Range("N96").Value = -7403.5
$ws.Range("H105").Value = 2043.875
$ws.Range("I105").Value = 2121.7144
$ws.Range("K105").Value = 2121.7144
$ws.Range("M105").Value = -374.7143999999998
$ws.Range("H122").Value = 782.36365
$ws.Range("I122").Value = 839.5
$ws.Range("J122").Value = 713.8
$ws.Range("K122").Value = 2518.5
$ws.Range("L122").Value = 2141.4
$ws.Range("M122").Value = -68.5
$ws.Range("N122").Value = -7041.4
$ws.Range("H132").Value = 1841.5333
$ws.Range("I132").Value = 1846.3462
$ws.Range("K132").Value = 5539.0386
$ws.Range("M132").Value = -3009.0386
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H47").Value = 100
$ws.Range("I47").Value = 100
$ws.Range("K47").Value = 100
$ws.Range("M47").Value = 466
$ws.Range("H141").Value = 27999
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Worksheet: CUL (sheet5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31783546
$ws.Range("J4").Value = 73456.78999999999
$ws.Range("L4").Value = 220370.37
$ws.Range("N4").Value = -220594.37
$ws.Range("H127").Value = 57490
$ws.Range("J127").Value = 57490
$ws.Range("L127").Value = 172470
$ws.Range("N127").Value = -182390
$ws.Range("H134").Value = 4678.3335
$ws.Range("J134").Value = 5495
$ws.Range("L134").Value = 16485
$ws.Range("N134").Value = -26625
$ws.Range("H140").Value = 1769.5483
$ws.Range("I140").Value = 1169.381
$ws.Range("K140").Value = 3508.143
$ws.Range("M140").Value = 1671.857
$ws.Range("H141").Value = 3506.2307
$ws.Range("I141").Value = 3417.2727
$ws.Range("K141").Value = 10251.8181
$ws.Range("M141").Value = -5071.8181
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H113").Value = 111114340
$ws.Range("I113").Value = 200
$ws.Range("K113").Value = 600
$ws.Range("M113").Value = 1570

# --- Worksheet: GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 396745.16
$ws.Range("I102").Value = 458836.25
$ws.Range("J102").Value = 6458.2856
$ws.Range("K102").Value = 458836.25
$ws.Range("L102").Value = 6458.2856
$ws.Range("M102").Value = -457214.25
$ws.Range("N102").Value = -9702.285599999999
$ws.Range("H126").Value = 20837382
$ws.Range("I126").Value = 50003180
$ws.Range("J126").Value = 4670.2144
$ws.Range("K126").Value = 150009540
$ws.Range("L126").Value = 14010.6432
$ws.Range("M126").Value = -150007070
$ws.Range("N126").Value = -18950.6432

# --- Worksheet: LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3605.0588
$ws.Range("I7").Value = 3552.0667
$ws.Range("J7").Value = 4002.5
$ws.Range("K7").Value = 3552.0667
$ws.Range("L7").Value = 4002.5
$ws.Range("M7").Value = -3440.0667
$ws.Range("N7").Value = -4226.5
$ws.Range("H11").Value = 13999.4
$ws.Range("J11").Value = 13999.4
$ws.Range("L11").Value = 13999.4
$ws.Range("N11").Value = -14279.4
$ws.Range("H22").Value = 3290
$ws.Range("I22").Value = 3288.3333
$ws.Range("J22").Value = 3292
$ws.Range("K22").Value = 3288.3333
$ws.Range("L22").Value = 3292
$ws.Range("M22").Value = -2993.3333
$ws.Range("N22").Value = -3882
$ws.Range("H27").Value = 3290
$ws.Range("I27").Value = 3288.3333
$ws.Range("J27").Value = 3292
$ws.Range("K27").Value = 3288.3333
$ws.Range("L27").Value = 3292
$ws.Range("M27").Value = -3181.3333
$ws.Range("N27").Value = -3506
$ws.Range("H122").Value = 4560.0977
$ws.Range("J122").Value = 6310.5884
$ws.Range("L122").Value = 18931.7652
$ws.Range("N122").Value = -23831.7652
$ws.Range("H126").Value = 3605.0588
$ws.Range("I126").Value = 3552.0667
$ws.Range("J126").Value = 4002.5
$ws.Range("K126").Value = 10656.2001
$ws.Range("L126").Value = 12007.5
$ws.Range("M126").Value = -8186.2001
$ws.Range("N126").Value = -16947.5
$ws.Range("H136").Value = 3391758.5
$ws.Range("I136").Value = 1875.2979
$ws.Range("J136").Value = 16668801
$ws.Range("K136").Value = 5625.893700000001
$ws.Range("L136").Value = 50006403
$ws.Range("M136").Value = -3075.893700000001
$ws.Range("N136").Value = -50011503

# --- Worksheet: WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 400
$ws.Range("K17").Value = 400
$ws.Range("M17").Value = -228
$ws.Range("H22").Value = 10503.167
$ws.Range("J22").Value = 12205.8
$ws.Range("L22").Value = 12205.8
$ws.Range("N22").Value = -12791.8
$ws.Range("H40").Value = 29439.572
$ws.Range("I40").Value = 28664.666
$ws.Range("J40").Value = 30020.75
$ws.Range("K40").Value = 28664.666
$ws.Range("L40").Value = 30020.75
$ws.Range("M40").Value = -28515.666
$ws.Range("N40").Value = -30318.75
$ws.Range("H81").Value = 1808.2667
$ws.Range("I81").Value = 1808.2667
$ws.Range("K81").Value = 3616.5334
$ws.Range("M81").Value = -2555.5334
$ws.Range("H84").Value = 1808.2667
$ws.Range("I84").Value = 1808.2667
$ws.Range("K84").Value = 18082.667
$ws.Range("M84").Value = -12778.667
$ws.Range("H136").Value = 2901.2954
$ws.Range("I136").Value = 2690.8108
$ws.Range("J136").Value = 4013.8572
$ws.Range("K136").Value = 8072.432400000001
$ws.Range("L136").Value = 12041.5716
$ws.Range("M136").Value = -5522.432400000001
$ws.Range("N136").Value = -17141.5716
